$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.197179061263853
$ws.Range("C2").Value = 0.3542999774293492
$ws.Range("D2").Value = 0.07977369863188244
$ws.Range("E2").Value = 0.4262032743756521
$ws.Range("G2").Value = 0.2241307714881842
$ws.Range("H2").Value = 0.3804231017443058
$ws.Range("I2").Value = 0.2288158210401221
$ws.Range("O2").Value = 1.12163897171763

$ws.Range("B3").Value = 1.046789645594743
$ws.Range("C3").Value = 0.3111477231075526
$ws.Range("D3").Value = 0.07219123885781187
$ws.Range("E3").Value = 0.3717201788608548
$ws.Range("G3").Value = 0.2216185945061326
$ws.Range("H3").Value = 0.3842275100089552
$ws.Range("I3").Value = 0.2360415229776045
$ws.Range("O3").Value = 1.123950320604308

$ws.Range("B4").Value = 0.9541225562241493
$ws.Range("C4").Value = 0.2845318556371979
$ws.Range("D4").Value = 0.06757029860568764
$ws.Range("E4").Value = 0.3383604125574067
$ws.Range("G4").Value = 0.220431523554403
$ws.Range("H4").Value = 0.3868910134436589
$ws.Range("I4").Value = 0.2407968291570963
$ws.Range("O4").Value = 1.126797602685457

$ws.Range("B5").Value = 0.9162799075818953
$ws.Range("C5").Value = 0.2736560495529545
$ws.Range("D5").Value = 0.06569593246118188
$ws.Range("E5").Value = 0.3247875541065213
$ws.Range("G5").Value = 0.2200364023568397
$ws.Range("H5").Value = 0.3880585369501546
$ws.Range("I5").Value = 0.2428143937961043
$ws.Range("O5").Value = 1.128314878432377

$ws.Range("B6").Value = 0.9099913966365989
$ws.Range("C6").Value = 0.2718483589634388
$ws.Range("D6").Value = 0.06538522037196515
$ws.Range("E6").Value = 0.3225350359862063
$ws.Range("G6").Value = 0.2199761231683866
$ws.Range("H6").Value = 0.3882573561641038
$ws.Range("I6").Value = 0.2431542131780837
$ws.Range("O6").Value = 1.128588319770074

$ws.Range("B7").Value = 0.9536125181534203
$ws.Range("C7").Value = 0.2843852998506122
$ws.Range("D7").Value = 0.06754498503823925
$ws.Range("E7").Value = 0.3381772795914202
$ws.Range("G7").Value = 0.2204258369918648
$ws.Range("H7").Value = 0.3869064268886362
$ws.Range("I7").Value = 0.2408237164121285
$ws.Range("O7").Value = 1.126816622599833

$ws.Range("B8").Value = 1.145394042866371
$ws.Range("C8").Value = 0.3394463761880786
$ws.Range("D8").Value = 0.07715204189865688
$ws.Range("E8").Value = 0.4073970310551118
$ws.Range("G8").Value = 0.2231903522899756
$ws.Range("H8").Value = 0.3816667222062335
$ws.Range("I8").Value = 0.2312408581325869
$ws.Range("O8").Value = 1.122138240673891

$ws.Range("B9").Value = 1.518798957898639
$ws.Range("C9").Value = 0.446446057934736
$ws.Range("D9").Value = 0.09626915846183692
$ws.Range("E9").Value = 0.5439668457514131
$ws.Range("G9").Value = 0.2314661566143172
$ws.Range("H9").Value = 0.3740022780989278
$ws.Range("I9").Value = 0.2149952446112753
$ws.Range("O9").Value = 1.124392713933844

$ws.Range("B10").Value = 1.791426862770265
$ws.Range("C10").Value = 0.524443957293272
$ws.Range("D10").Value = 0.1104883493353128
$ws.Range("E10").Value = 0.6449468702703882
$ws.Range("G10").Value = 0.2393352877735708
$ws.Range("H10").Value = 0.3699787154206149
$ws.Range("I10").Value = 0.2046357307084943
$ws.Range("O10").Value = 1.133153516930889

$ws.Range("B11").Value = 1.91506569944994
$ws.Range("C11").Value = 0.5597900229861921
$ws.Range("D11").Value = 0.1169958458894911
$ws.Range("E11").Value = 0.6910548711811515
$ws.Range("G11").Value = 0.2433141192691011
$ws.Range("H11").Value = 0.3685008250765947
$ws.Range("I11").Value = 0.2002702662215157
$ws.Range("O11").Value = 1.138711377249251

$ws.Range("B12").Value = 1.96182793499662
$ws.Range("C12").Value = 0.5731546940993439
$ws.Range("D12").Value = 0.1194657439208555
$ws.Range("E12").Value = 0.7085417333968564
$ws.Range("G12").Value = 0.2448790261817635
$ws.Range("H12").Value = 0.3679921504560184
$ws.Range("I12").Value = 0.1986675267657922
$ws.Range("O12").Value = 1.141044487386807

$ws.Range("B13").Value = 1.951759427485456
$ws.Range("C13").Value = 0.5702772775582616
$ws.Range("D13").Value = 0.1189335554786197
$ws.Range("E13").Value = 0.7047744011567971
$ws.Range("G13").Value = 0.2445393935775257
$ws.Range("H13").Value = 0.3680994311071686
$ws.Range("I13").Value = 0.1990104581097167
$ws.Range("O13").Value = 1.140531811824218

$ws.Range("B14").Value = 1.918914014456504
$ws.Range("C14").Value = 0.5608899492234514
$ws.Range("D14").Value = 0.1171989324983969
$ws.Range("E14").Value = 0.6924929775180004
$ws.Range("G14").Value = 0.2434416938142618
$ws.Range("H14").Value = 0.3684579529299157
$ws.Range("I14").Value = 0.200137395865978
$ws.Range("O14").Value = 1.138898731339623

$ws.Range("B15").Value = 1.898787744899607
$ws.Range("C15").Value = 0.5551373014575915
$ws.Range("D15").Value = 0.1161371624392871
$ws.Range("E15").Value = 0.6849737982443997
$ws.Range("G15").Value = 0.2427769258818415
$ws.Range("H15").Value = 0.3686842043312026
$ws.Range("I15").Value = 0.2008342504560066
$ws.Range("O15").Value = 1.137928245358836

$ws.Range("B16").Value = 1.783338908396104
$ws.Range("C16").Value = 0.5221312273677086
$ws.Range("D16").Value = 0.1100638584119054
$ws.Range("E16").Value = 0.6419372355227608
$ws.Range("G16").Value = 0.2390833629324334
$ws.Range("H16").Value = 0.3700824158534886
$ws.Range("I16").Value = 0.2049280477015234
$ws.Range("O16").Value = 1.132822141038275

$ws.Range("B17").Value = 1.712415546329396
$ws.Range("C17").Value = 0.5018479016338233
$ws.Range("D17").Value = 0.1063481196005966
$ws.Range("E17").Value = 0.6155811392726633
$ws.Range("G17").Value = 0.236920287510614
$ws.Range("H17").Value = 0.3710306489859079
$ws.Range("I17").Value = 0.2075286882477236
$ws.Range("O17").Value = 1.130094138349079

$ws.Range("B18").Value = 1.671586583016278
$ws.Range("C18").Value = 0.4901687452544365
$ws.Range("D18").Value = 0.104214605091201
$ws.Range("E18").Value = 0.6004378438414477
$ws.Range("G18").Value = 0.2357136623678855
$ws.Range("H18").Value = 0.371609204397231
$ws.Range("I18").Value = 0.2090571544803481
$ws.Range("O18").Value = 1.128672936836324

$ws.Range("B19").Value = 1.657756533851796
$ws.Range("C19").Value = 0.4862122173449279
$ws.Range("D19").Value = 0.1034928650341271
$ws.Range("E19").Value = 0.5953132751862995
$ws.Range("G19").Value = 0.2353115388318372
$ws.Range("H19").Value = 0.3718107792523
$ws.Range("I19").Value = 0.2095802602612249
$ws.Range("O19").Value = 1.128217068557518

$ws.Range("B20").Value = 1.719969173279537
$ws.Range("C20").Value = 0.5040084176996515
$ws.Range("D20").Value = 0.1067432852717332
$ws.Range("E20").Value = 0.6183851163115719
$ws.Range("G20").Value = 0.2371466616339291
$ws.Range("H20").Value = 0.3709262742551402
$ws.Range("I20").Value = 0.2072484635625997
$ws.Range("O20").Value = 1.130369217582228

$ws.Range("B21").Value = 1.928563076729006
$ws.Range("C21").Value = 0.5636477859827096
$ws.Range("D21").Value = 0.1177082797625957
$ws.Range("E21").Value = 0.696099586352986
$ws.Range("G21").Value = 0.243762528537971
$ws.Range("H21").Value = 0.368351260752604
$ws.Range("I21").Value = 0.1998050163852483
$ws.Range("O21").Value = 1.139372187399886

$ws.Range("B22").Value = 2.064557322987525
$ws.Range("C22").Value = 0.6025080224183625
$ws.Range("D22").Value = 0.1249074839498547
$ws.Range("E22").Value = 0.7470474599680585
$ws.Range("G22").Value = 0.2484260504959792
$ws.Range("H22").Value = 0.3669655696257053
$ws.Range("I22").Value = 0.1952340815465607
$ws.Range("O22").Value = 1.146588930451401

$ws.Range("B23").Value = 1.99200594861037
$ws.Range("C23").Value = 0.5817785407214728
$ws.Range("D23").Value = 0.1210621100664753
$ws.Range("E23").Value = 0.7198405658624409
$ws.Range("G23").Value = 0.245905691300834
$ws.Range("H23").Value = 0.3676778469397419
$ws.Range("I23").Value = 0.1976466448584233
$ws.Range("O23").Value = 1.1426144946129

$ws.Range("B24").Value = 1.71655434712244
$ws.Range("C24").Value = 0.5030317044207777
$ws.Range("D24").Value = 0.1065646224243437
$ws.Range("E24").Value = 0.6171174099114012
$ws.Range("G24").Value = 0.237044202838419
$ws.Range("H24").Value = 0.3709733580370767
$ws.Range("I24").Value = 0.207375049373006
$ws.Range("O24").Value = 1.130244396019094

$ws.Range("B25").Value = 1.418078148630229
$ws.Range("C25").Value = 0.4176063291037053
$ws.Range("D25").Value = 0.09106727232408218
$ws.Range("E25").Value = 0.5069182298679209
$ws.Range("G25").Value = 0.228916700160795
$ws.Range("H25").Value = 0.3757945785605585
$ws.Range("I25").Value = 0.2191150508832322
$ws.Range("O25").Value = 1.122545123907742
